# Apply "output generated at 456a3b4" data refresh to 上海-漫展信息.xlsx
# Sheets: 1=展览(Exhibitions) 2=演出(Shows) 3=本地生活(Local life) 4=全部类型(All types)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------
# Sheet "展览" (1): refreshed "想去人数" (F) counts
# ---------------------------------------------------------------
$ws1.Range("F3").Value = 3326
$ws1.Range("F6").Value = 7762
$ws1.Range("F9").Value = 1155
$ws1.Range("F10").Value = 1063
$ws1.Range("F14").Value = 1770
$ws1.Range("F15").Value = 375
$ws1.Range("F17").Value = 2374
$ws1.Range("F18").Value = 1024
$ws1.Range("F20").Value = 1033
$ws1.Range("F21").Value = 1055
$ws1.Range("F22").Value = 6370
$ws1.Range("F23").Value = 7014
$ws1.Range("F24").Value = 415
$ws1.Range("F26").Value = 1094
$ws1.Range("F29").Value = 524
$ws1.Range("F30").Value = 1082
$ws1.Range("F31").Value = 1043
$ws1.Range("F32").Value = 526
$ws1.Range("F33").Value = 526
$ws1.Range("F35").Value = 91
$ws1.Range("F39").Value = 422
$ws1.Range("F40").Value = 345
$ws1.Range("F41").Value = 1277
$ws1.Range("F42").Value = 3275
$ws1.Range("F43").Value = 615
$ws1.Range("F44").Value = 721
$ws1.Range("F45").Value = 484
$ws1.Range("F47").Value = 108
$ws1.Range("F48").Value = 100
$ws1.Range("F49").Value = 487
$ws1.Range("F50").Value = 68

# ---------------------------------------------------------------
# Sheet "演出" (2): refreshed F (想去人数) and G (最低票价) values
# ---------------------------------------------------------------
$ws2.Range("F4").Value = 384
$ws2.Range("F5").Value = 649
$ws2.Range("F9").Value = 83
$ws2.Range("F14").Value = 3
$ws2.Range("F24").Value = 5
$ws2.Range("F26").Value = 13
$ws2.Range("F27").Value = 6629
$ws2.Range("F33").Value = 2

$ws2.Range("G17").Value = 380
$ws2.Range("G27").Value = 680

# ---------------------------------------------------------------
# Sheet "本地生活" (3): refreshed F values + G8 resolves from
# "已售罄" (sold out, text) to a real minimum price (number), plus
# one brand-new row (12) for a newly-listed event.
# ---------------------------------------------------------------
$ws3.Range("F4").Value = 2010
$ws3.Range("F5").Value = 1330
$ws3.Range("F7").Value = 563
$ws3.Range("F8").Value = 2160
$ws3.Range("G8").Value = 10
$ws3.Range("F9").Value = 8960
$ws3.Range("F10").Value = 1102
$ws3.Range("F11").Value = 92

$ws3.Range("A12").Value = 11
$ws3.Range("A11").Copy()
$ws3.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Range("B12").Value = "'2024-09-09"
$ws3.Range("C12").Value = "上海·日漫咖啡体验"
$ws3.Range("D12").Value = "虹桥路1438号高岛屋百货6楼 Oasis漫画喫茶"
$ws3.Range("E12").Value = "2024.09.09 10:00-12.31 22:00"
$ws3.Range("F12").Value = 0
$ws3.Range("G12").Value = 60
$ws3.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=91993"
$ws3.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202409/IV5rInWT1725347808557.jpeg"

# ---------------------------------------------------------------
# Sheet "全部类型" (4): refreshed F values, plus the curated
# rows 4-6 are re-pointed: the "THE哆啦A梦展" entry drops off the
# list, rows 5 & 6 slide up into 4 & 5 (with refreshed numbers),
# and row 6 is replaced by a brand-new entry (东方明珠 pop-up).
# ---------------------------------------------------------------
$ws4.Range("F2").Value = 3326
$ws4.Range("F3").Value = 2010

$ws4.Range("C4").Value = "上海·艺术与潮流·遇见EVA 中国首展"
$ws4.Range("D4").Value = "西藏北路166号 静安大悦城北座"
$ws4.Range("E4").Value = "2024.08.02 10:00-10.07 22:00"
$ws4.Range("F4").Value = 7762
$ws4.Range("G4").Value = 89
$ws4.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=89161"
$ws4.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202407/z8YTdxA71720679877329.jpeg"

$ws4.Range("B5").Value = "'2024-08-05"
$ws4.Range("C5").Value = "上海·名侦探柯南 连载30周年纪念展"
$ws4.Range("D5").Value = "南京西路2-68号 新世界城11楼"
$ws4.Range("E5").Value = "2024.08.05 10:00-10.07 22:00"
$ws4.Range("F5").Value = 1330
$ws4.Range("G5").Value = 109
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=89870"
$ws4.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202407/35thNBrO1721035918311.png"

$ws4.Range("B6").Value = "'2024-08-17"
$ws4.Range("C6").Value = "上海·东方明珠·「光与夜之恋 × 线条小狗 ×爱胖达文化 」线条大作战主题店"
$ws4.Range("D6").Value = "世纪大道1号 东方明珠电视塔城市广场商场"
$ws4.Range("E6").Value = "2024.08.17 00:00-10.27 23:59"
$ws4.Range("F6").Value = 2160
$ws4.Range("G6").Value = 10
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=90444"
$ws4.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202408/qUE9n4UR1723020534077.png"

$ws4.Range("F8").Value = 1102
$ws4.Range("F9").Value = 92
$ws4.Range("F10").Value = 1155
$ws4.Range("F11").Value = 1063
$ws4.Range("F15").Value = 375
$ws4.Range("F18").Value = 2374
$ws4.Range("F19").Value = 1024
$ws4.Range("F20").Value = 1033
$ws4.Range("F21").Value = 1055
$ws4.Range("F22").Value = 6370
$ws4.Range("F23").Value = 7014
$ws4.Range("F24").Value = 415
$ws4.Range("F26").Value = 1094
$ws4.Range("F29").Value = 524
$ws4.Range("F30").Value = 1043
$ws4.Range("F31").Value = 526
$ws4.Range("F33").Value = 91
$ws4.Range("F37").Value = 422
$ws4.Range("F38").Value = 345
$ws4.Range("F40").Value = 3275
$ws4.Range("F41").Value = 615
$ws4.Range("F42").Value = 721
$ws4.Range("F43").Value = 484
$ws4.Range("F44").Value = 108
$ws4.Range("F45").Value = 100
$ws4.Range("F48").Value = 68
